$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Row 8 content (Case4_12v_3f) - replace the old "data format" warning with
#    a full set of populated cells across the row.
# ---------------------------------------------------------------------------
$ws.Range("C8").Value = "simulated dataset"
$ws.Range("D8").Value = "Not available in source"
$ws.Range("E8").Value = "completely standardized,  uses value command to fix unstandardized loading to 1, Maximum likelihood"
$ws.Range("F8").Value = "chi2, SRMR, RMSEA, TLI: same range; CFI: identical"
$ws.Range("G8").Value = "not available "
$ws.Range("H8").Value = "not available "
$ws.Range("I8").Value = "same range"
$ws.Range("J8").Value = "Not available "
$ws.Range("K8").Value = "x4-social and x1-x2: identified as modification indices; x11-x12 and x9-x10: not identified"
$ws.Range("L8").Value = "not available in source"
$ws.Range("M8").Value = """fit indices are consistent with good model fit"""
$ws.Range("N8").Value = 'Source: "standardized residuals (5.04) and modification indices (d12,11 = 25.94) indicate  that the relationship between these items has not been adequately reproduced by the  model' + [char]0x2019 + 's parameter estimates."' + [char]10 + ' '
$ws.Range("P8").Value = "Source: ""All freely estimated parameters are statistically significant."""

# ---------------------------------------------------------------------------
# 2) Row 9 content (Case5_8v_2f) - fill in the previously empty cells.
# ---------------------------------------------------------------------------
$ws.Range("D9").Value = "Not available in source"
$ws.Range("E9").Value = "first indicator as marker variable (factor loading of first indicator=1), not standardized"
$ws.Range("J9").Value = "identical"
$ws.Range("K9").Value = "same modification indices identified, but different values"
$ws.Range("L9").Value = "Values of the CFI and SRMR are, respectively, .959 and .072, and neither result is  clearly problematic."
$ws.Range("M9").Value = "-"
$ws.Range("N9").Value = 'Source:"Most of the larger and positive residuals are between Hand Movements and other tasks specified to measure the other factor. Because the standardized  pattern coefficient of Hand Movements is at least moderate (.497; Table 13.3), it is possible that this task may measure both factors."'
$ws.Range("O9").Value = "not available in source"
$ws.Range("P9").Value = "Source: ""handmov is not specified to measure Simul"""

# ---------------------------------------------------------------------------
# 3) Re-apply formatting to row 8 / row 9 cells using existing cells in the
#    sheet as format donors (keeps the style table de-duplicated, exactly as
#    Excel itself would when the same combination of fill/border/alignment is
#    reused elsewhere in the workbook).
# ---------------------------------------------------------------------------

# style "6" (no fill, thin border, wrap, top-aligned) - default data cell
$ws.Range("D6").Copy()
$ws.Range("D8,G8,H8,J8,L8,D9,E9,M9,O9").PasteSpecial(-4122)

# style "7" (green fill, thin border, wrap, top-aligned)
$ws.Range("C6").Copy()
$ws.Range("C8,F8,M8,G9,J9,K9,L9").PasteSpecial(-4122)

# style "8" (purple fill, thin border, wrap, top-aligned)
$ws.Range("P6").Copy()
$ws.Range("P8,N9,P9").PasteSpecial(-4122)

# style "10" (orange fill, thin border, wrap, top-aligned)
$ws.Range("I7").Copy()
$ws.Range("I8,K8,I9").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 4) N8:O8 are merged and get a brand-new style: purple fill, centered,
#    top-aligned, wrapped text, with the shared inner border suppressed
#    (N8 keeps its left/top/bottom border, O8 keeps its right/top/bottom
#    border - matching how every other merged header cell in row 2 is built).
# ---------------------------------------------------------------------------
$mergeRange = $ws.Range("N8:O8")
$mergeRange.Interior.Color = 10498160
$mergeRange.HorizontalAlignment = -4108
$mergeRange.VerticalAlignment = -4160
$mergeRange.WrapText = $true

foreach ($addr in @("N8", "O8")) {
    $cell = $ws.Range($addr)
    $cell.Borders.Item(8).LineStyle = 1
    $cell.Borders.Item(8).Weight = 2
    $cell.Borders.Item(9).LineStyle = 1
    $cell.Borders.Item(9).Weight = 2
}
$ws.Range("N8").Borders.Item(7).LineStyle = 1
$ws.Range("N8").Borders.Item(7).Weight = 2
$ws.Range("O8").Borders.Item(10).LineStyle = 1
$ws.Range("O8").Borders.Item(10).Weight = 2

$ws.Range("N8:O8").Merge()

# ---------------------------------------------------------------------------
# 5) Freeze panes (2 columns / 4 rows) and restore the per-pane selections.
# ---------------------------------------------------------------------------
$ws.Range("C5").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("E9").Select()
